# Add new benchmark row (v11-18000) to Sheet1 describing the fixed
# "reversing direction" bug, matching the author's "V11" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$row = 13

$ws.Cells.Item($row, 1).Value = "v11-18000"
$ws.Cells.Item($row, 2).Value = 64
$ws.Cells.Item($row, 3).Value = "Fixed severe bug when reversing direction"
$ws.Cells.Item($row, 4).Value = 85
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 1838.94
$ws.Cells.Item($row, 10).Value = 1410
$ws.Cells.Item($row, 11).Value = 2860
$ws.Cells.Item($row, 12).Value = 171.71
$ws.Cells.Item($row, 13).Value = 133
$ws.Cells.Item($row, 14).Value = 230
$ws.Cells.Item($row, 15).Value = 0.11
$ws.Cells.Item($row, 16).Value = 0
$ws.Cells.Item($row, 17).Value = 2
$ws.Cells.Item($row, 18).Value = 95.27
$ws.Cells.Item($row, 19).Value = 58.7
$ws.Cells.Item($row, 20).Value = 223

$ws.Range("C14").Select()
